$wb = $excel.ActiveWorkbook

# --- Rename sheets (refresh timestamps in generated task-order names) ---
$wb.Worksheets.Item("GNG_TO-16511687280351787").Name  = "GNG_TO-1651255553244722"
$wb.Worksheets.Item("NB_TO-16511687314070442").Name    = "NB_TO-16512555558837872"
$wb.Worksheets.Item("RS_TO-1651168731408046").Name     = "RS_TO-165125555588579"
$wb.Worksheets.Item("TOL_TO-16511687314548554").Name   = "TOL_TO-16512555559497843"
$wb.Worksheets.Item("vSAT_TO-1651168731516651").Name   = "vSAT_TO-16512555560277808"

# --- GNG_TO sheet: update stim file names ---
$ws1 = $wb.Worksheets.Item("GNG_TO-1651255553244722")
$ws1.Range("B2").Value = "go_stims-16512555532097223.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555532267325.csv"
$ws1.Range("B4").Value = "go_stims-16512555532287216.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555532427206.csv"

# --- NB_TO sheet: update stim file names ---
$ws2 = $wb.Worksheets.Item("NB_TO-16512555558837872")
$ws2.Range("B2").Value  = "ZB-match_9-16512555534386263.csv"
$ws2.Range("B3").Value  = "ZB-match_1-16512555532557232.csv"
$ws2.Range("B4").Value  = "TB-16512555558687804.csv"
$ws2.Range("B5").Value  = "ZB-match_5-16512555533007185.csv"
$ws2.Range("B6").Value  = "OB-16512555541216266.csv"
$ws2.Range("B7").Value  = "TB-16512555558137946.csv"
$ws2.Range("B8").Value  = "OB-16512555537806268.csv"
$ws2.Range("B9").Value  = "OB-16512555535056272.csv"
$ws2.Range("B10").Value = "TB-16512555550676258.csv"

# --- TOL_TO sheet: update stim file names ---
$ws4 = $wb.Worksheets.Item("TOL_TO-16512555559497843")
$ws4.Range("B2").Value = "MM_stims-16512555559157825.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555558937893.csv"
$ws4.Range("B4").Value = "MM_stims-165125555593178.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255555916784.csv"
$ws4.Range("B6").Value = "MM_stims-165125555594778.csv"
$ws4.Range("B7").Value = "ZM_stims-165125555593278.csv"

# --- vSAT_TO sheet: update stim file names ---
$ws5 = $wb.Worksheets.Item("vSAT_TO-16512555560277808")
$ws5.Range("B2").Value = "SAT_stims-16512555559537828.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555559807801.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555559957795.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555560117812.csv"
